$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new record (2026/01/28 23:00, ranking 26) was logged right after the
# existing 2026/01/28 19:00 entry (row 712) and before the 2026/12/29
# 13:00 entry (row 713). Insert a new row at 713 so the remaining rows
# shift down by one, then populate the new row with this record.
$ws.Rows("713:713").Insert()

# Column A stores the date as plain text (matching the rest of the
# column), not an Excel date serial, so force a text number format
# before assigning, then restore the default "Normal" style so no
# stray formatting is left on the cell.
$ws.Cells.Item(713, 1).NumberFormat = "@"
$ws.Cells.Item(713, 1).Value = "2026/01/28"
$ws.Cells.Item(713, 1).Style = "Normal"

$ws.Cells.Item(713, 2).Value = "水"
$ws.Cells.Item(713, 3).Value = 23
$ws.Cells.Item(713, 4).Value = 26
